# Fix the error of material in 10
# The demo dates in the "日期条件格式" and "文字条件格式" sheets were stale
# (from 2014); bump every date value forward by exactly one year (365 days)
# so the conditional-formatting examples show "current" looking data again.

$wb = $excel.ActiveWorkbook

# --- Sheet 2: 日期条件格式 -------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()

$ws2.Range("B2").Value  = 42193
$ws2.Range("B3").Value  = 42228
$ws2.Range("B4").Value  = 42200
$ws2.Range("B5").Value  = 42212
$ws2.Range("B6").Value  = 42220
$ws2.Range("B7").Value  = 42208
$ws2.Range("B8").Value  = 42209
$ws2.Range("B9").Value  = 42220
$ws2.Range("B10").Value = 42230

$ws2.Range("B10").Select()

# --- Sheet 3: 文字条件格式 -------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()

$ws3.Range("D2").Value  = 42150
$ws3.Range("D3").Value  = 42157
$ws3.Range("D4").Value  = 42037
$ws3.Range("D5").Value  = 42197
$ws3.Range("D6").Value  = 42091
$ws3.Range("D7").Value  = 42122
$ws3.Range("D8").Value  = 42157
$ws3.Range("D9").Value  = 42091
$ws3.Range("D10").Value = 42163

$ws3.Range("D10").Select()

# --- Leave the workbook back on the first sheet (条件格式), as in the source ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
